# Adding test cases for watch list (TestCase_E35, TestCase_E36)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 36 - copy formatting from row 35 (same style pattern: fill-highlighted description cell)
$ws.Range("A35:E35").Copy($ws.Range("A36:E36"))

# Row 37 - copy formatting from row 33 (plain description cell, no fill)
$ws.Range("A33:E33").Copy($ws.Range("A37:E37"))

# Fill in row 36 values (Description set first, then TCID, Jira id, to match shared-string order)
$ws.Range("C36").Value = "Verify that user is able to convert his public watchlist to private"
$ws.Range("A36").Value = "TestCase_E35"
$ws.Range("B36").Value = "OPQA-330"
$ws.Range("D36").Value = "Y"
$ws.Range("E36").Value = "PASS"

# Fill in row 37 values
$ws.Range("C37").Value = "Verify that user is able to see the watchlist items by content type"
$ws.Range("A37").Value = "TestCase_E36"
$ws.Range("B37").Value = "OPQA-618"
$ws.Range("D37").Value = "Y"
$ws.Range("E37").Value = "PASS"

# Update the sheet view: drop the frozen/topLeft C1 view and select A3 instead
$ws.Range("A3").Select()

# Best-effort: reflect the new workbook window size recorded in the saved file
$win = $wb.Windows.Item(1)
$win.Width = 14310
$win.Height = 10125
